$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.502.53"
$ws.Range("E2").Value = "  -2.88%  "
$ws.Range("D3").Value = "1.801.75"
$ws.Range("E3").Value = "  -2.27%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").Value = "'229.17"
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("E6").Value = "  -1.41%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").Value = "'39.33"
$ws.Range("E8").Value = "  -11.44%  "
$ws.Range("E9").Value = "  +2.90%  "
$ws.Range("D10").Value = "'0.0678"
$ws.Range("E10").Value = "  -2.81%  "
$ws.Range("E11").Value = "  -2.16%  "
$ws.Range("D12").Value = "2.061.32"
$ws.Range("E12").Value = "  -2.34%  "
$ws.Range("D13").Value = "'11.09"
$ws.Range("E13").Value = "  -1.95%  "
$ws.Range("E14").Value = "  -2.39%  "
$ws.Range("D15").Value = "1.790.29"
$ws.Range("E15").Value = "  -2.94%  "
$ws.Range("D16").Value = "'4.56"
$ws.Range("E16").Value = "  -3.56%  "
$ws.Range("D17").Value = "34.362.70"
$ws.Range("E17").Value = "  -3.23%  "
$ws.Range("D18").Value = "'69.03"
$ws.Range("E18").Value = "  -2.05%  "
$ws.Range("E19").Value = "  -2.96%  "
$ws.Range("D20").Value = "'239.35"
$ws.Range("E20").Value = "  -1.98%  "
$ws.Range("D21").Value = "'11.77"
$ws.Range("E21").Value = "  -2.76%  "
$ws.Range("D22").Value = "'4.68"
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("E24").Value = "  -2.47%  "
$ws.Range("D25").Value = "'172.88"
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("D26").Value = "'7.73"
$ws.Range("E26").Value = "  -3.39%  "
$ws.Range("D27").Value = "'17.17"
$ws.Range("E27").Value = "  -3.82%  "
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("E31").Value = "  +1.04%  "
$ws.Range("E32").Value = "  -2.40%  "
$ws.Range("E33").Value = "  -5.28%  "
$ws.Range("E34").Value = "  +7.07%  "
$ws.Range("D35").Value = "'1.78"
$ws.Range("E35").Value = "  -3.23%  "
$ws.Range("D36").Value = "'0.692"
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("D37").Value = "'90.33"
$ws.Range("E37").Value = "  -5.86%  "
$ws.Range("E38").Value = "  +5.06%  "
$ws.Range("D39").Value = "1.323.26"
$ws.Range("E39").Value = "  -1.73%  "
$ws.Range("E40").Value = "  -3.12%  "
$ws.Range("E41").Value = "  -5.73%  "
$ws.Range("D42").Value = "'14.21"
$ws.Range("E42").Value = "  -6.62%  "
$ws.Range("D43").Value = "'2.39"
$ws.Range("E43").Value = "  -2.99%  "
$ws.Range("E44").Value = "  -9.32%  "
$ws.Range("E45").Value = "  -3.92%  "
$ws.Range("D46").Value = "'6.14"
$ws.Range("E46").Value = "  -2.11%  "
$ws.Range("E47").Value = "  -1.19%  "
$ws.Range("D48").Value = "1.986.07"
$ws.Range("E48").Value = "  -1.54%  "
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("E50").Value = "  +3.50%  "
$ws.Range("D51").Value = "'97.49"
$ws.Range("E51").Value = "  -4.84%  "
